$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28 (anchor G=27772)
$ws.Range("H28").Value = 904.8333
$ws.Range("J28").Value = 870
$ws.Range("L28").Value = 870
$ws.Range("N28").Value = -1840
# Row 41 (anchor G=5478)
$ws.Range("H41").Value = 449
$ws.Range("J41").Value = 473.5
$ws.Range("L41").Value = 473.5
$ws.Range("N41").Value = -1353.5
# Row 64 (anchor G=5506)
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 7000
$ws.Range("K64").Value = 7000
$ws.Range("M64").Value = -6752
# Row 67 (anchor G=5506)
$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 7000
$ws.Range("K67").Value = 7000
$ws.Range("M67").Value = -6142
# Row 70 (anchor G=12604)
$ws.Range("H70").Value = 2222
$ws.Range("J70").Value = 2500
$ws.Range("L70").Value = 7500
$ws.Range("N70").Value = -8040
# Row 73 (anchor G=12604)
$ws.Range("H73").Value = 2222
$ws.Range("J73").Value = 2500
$ws.Range("L73").Value = 7500
$ws.Range("N73").Value = -9372
# Row 86 (anchor G=12603)
$ws.Range("H86").Value = 7384
$ws.Range("J86").Value = 6076
$ws.Range("L86").Value = 6076
$ws.Range("N86").Value = -8322
# Row 89 (anchor G=12603)
$ws.Range("H89").Value = 7384
$ws.Range("J89").Value = 6076
$ws.Range("L89").Value = 30380
$ws.Range("N89").Value = -41612

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (anchor G=43999)
$ws.Range("H61").Value = 4054.625
$ws.Range("I61").Value = 3906.1667
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 3906.1667
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -3694.1667
$ws.Range("N61").Value = -4924
# Row 136 (anchor G=43999)
$ws.Range("H136").Value = 4054.625
$ws.Range("I136").Value = 3906.1667
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 11718.5001
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -9168.500100000001
$ws.Range("N136").Value = -18600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22 (anchor G=5092)
$ws.Range("H22").Value = 487.8
$ws.Range("I22").Value = 487.8
$ws.Range("K22").Value = 487.8
$ws.Range("M22").Value = -314.8
# Row 94 (anchor G=19939)
$ws.Range("H94").Value = 7564.625
$ws.Range("I94").Value = 4352.8335
$ws.Range("K94").Value = 4352.8335
$ws.Range("M94").Value = -3901.8335
# Row 105 (anchor G=19947)
$ws.Range("H105").Value = 8501.200000000001
$ws.Range("I105").Value = 10187.625
$ws.Range("J105").Value = 1755.5
$ws.Range("K105").Value = 10187.625
$ws.Range("L105").Value = 1755.5
$ws.Range("M105").Value = -8440.625
$ws.Range("N105").Value = -5249.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 122 (anchor G=36196)
$ws.Range("H122").Value = 1119.5
$ws.Range("I122").Value = 1255.1428
$ws.Range("K122").Value = 3765.4284
$ws.Range("M122").Value = -1315.4284
# Row 134 (anchor G=44020)
$ws.Range("H134").Value = 9750
$ws.Range("I134").Value = 9750
$ws.Range("K134").Value = 29250
$ws.Range("M134").Value = -26715

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 63 (anchor G=12866)
$ws.Range("H63").Value = 6500
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = -8251
$ws.Range("N63").Value = -31498
# Row 66 (anchor G=12866)
$ws.Range("H66").Value = 6500
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 27000
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -23256
$ws.Range("N66").Value = -97488
# Row 68 (anchor G=12895)
$ws.Range("H68").Value = 1360
$ws.Range("J68").Value = 1360
$ws.Range("L68").Value = 4080
$ws.Range("N68").Value = -5702
# Row 71 (anchor G=12895)
$ws.Range("H71").Value = 1360
$ws.Range("J71").Value = 1360
$ws.Range("L71").Value = 12240
$ws.Range("N71").Value = -20352
# Row 75 (anchor G=12863)
$ws.Range("H75").Value = 156.5
$ws.Range("I75").Value = 156.5
$ws.Range("K75").Value = 469.5
$ws.Range("M75").Value = 528.5
# Row 78 (anchor G=12863)
$ws.Range("H78").Value = 156.5
$ws.Range("I78").Value = 156.5
$ws.Range("K78").Value = 1408.5
$ws.Range("M78").Value = 3583.5
# Row 87 (anchor G=12864)
$ws.Range("H87").Value = 8249.5
$ws.Range("I87").Value = 15999
$ws.Range("J87").Value = 500
$ws.Range("K87").Value = 47997
$ws.Range("L87").Value = 1500
$ws.Range("M87").Value = -46749
$ws.Range("N87").Value = -3996
# Row 90 (anchor G=12864)
$ws.Range("H90").Value = 8249.5
$ws.Range("I90").Value = 15999
$ws.Range("J90").Value = 500
$ws.Range("K90").Value = 143991
$ws.Range("L90").Value = 4500
$ws.Range("M90").Value = -137751
$ws.Range("N90").Value = -16980
# Row 98 (anchor G=19843)
$ws.Range("H98").Value = 3299.25
$ws.Range("I98").Value = 4066
$ws.Range("J98").Value = 999
$ws.Range("K98").Value = 12198
$ws.Range("L98").Value = 2997
$ws.Range("M98").Value = -10700
$ws.Range("N98").Value = -5993
# Row 103 (anchor G=19839)
$ws.Range("H103").Value = 1328.5714
$ws.Range("I103").Value = 1045.75
$ws.Range("J103").Value = 1705.6666
$ws.Range("K103").Value = 3137.25
$ws.Range("L103").Value = 5116.9998
$ws.Range("M103").Value = -2258.25
$ws.Range("N103").Value = -6874.9998
# Row 114 (anchor G=27865)
$ws.Range("H114").Value = 1496.1428
$ws.Range("J114").Value = 1434.8
$ws.Range("L114").Value = 4304.4
$ws.Range("N114").Value = -10812.4
# Row 132 (anchor G=43972)
$ws.Range("H132").Value = 10330.5
$ws.Range("J132").Value = 13495.75
$ws.Range("L132").Value = 121461.75
$ws.Range("N132").Value = -126521.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97 (anchor G=19940)
$ws.Range("H97").Value = 7666
$ws.Range("I97").Value = 5999
$ws.Range("K97").Value = 5999
$ws.Range("M97").Value = -5503
# Row 113 (anchor G=27710)
$ws.Range("H113").Value = 2320.5
$ws.Range("I113").Value = 2702.5
$ws.Range("J113").Value = 2065.8333
$ws.Range("K113").Value = 2702.5
$ws.Range("L113").Value = 2065.8333
$ws.Range("M113").Value = -532.5
$ws.Range("N113").Value = -6405.8333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (anchor G=5277)
$ws.Range("H22").Value = 10243.538
$ws.Range("I22").Value = 10395.75
$ws.Range("K22").Value = 10395.75
$ws.Range("M22").Value = -10100.75
# Row 27 (anchor G=5277)
$ws.Range("H27").Value = 10243.538
$ws.Range("I27").Value = 10395.75
$ws.Range("K27").Value = 10395.75
$ws.Range("M27").Value = -10288.75
# Row 55 (anchor G=5284)
$ws.Range("H55").Value = 566.2857
$ws.Range("I55").Value = 648.5454999999999
$ws.Range("J55").Value = 264.66666
$ws.Range("K55").Value = 648.5454999999999
$ws.Range("L55").Value = 264.66666
$ws.Range("M55").Value = -475.5454999999999
$ws.Range("N55").Value = -610.66666
# Row 61 (anchor G=27740)
$ws.Range("H61").Value = 14573629
$ws.Range("I61").Value = 12752575
$ws.Range("J61").Value = 17001700
$ws.Range("K61").Value = 12752575
$ws.Range("L61").Value = 17001700
$ws.Range("M61").Value = -12752373
$ws.Range("N61").Value = -17002104
# Row 68 (anchor G=12563)
$ws.Range("H68").Value = 2109.5
$ws.Range("I68").Value = 2109.5
$ws.Range("K68").Value = 2109.5
$ws.Range("M68").Value = -1360.5
# Row 71 (anchor G=12563)
$ws.Range("H71").Value = 2109.5
$ws.Range("I71").Value = 2109.5
$ws.Range("K71").Value = 10547.5
$ws.Range("M71").Value = -6803.5
# Row 113 (anchor G=27740)
$ws.Range("H113").Value = 14573629
$ws.Range("I113").Value = 12752575
$ws.Range("J113").Value = 17001700
$ws.Range("K113").Value = 12752575
$ws.Range("L113").Value = 17001700
$ws.Range("M113").Value = -12750405
$ws.Range("N113").Value = -17006040
# Row 132 (anchor G=44058)
$ws.Range("H132").Value = 4999.5
$ws.Range("I132").Value = 3999.5
$ws.Range("J132").Value = 5499.5
$ws.Range("K132").Value = 11998.5
$ws.Range("L132").Value = 16498.5
$ws.Range("M132").Value = -9468.5
$ws.Range("N132").Value = -21558.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62 (anchor G=12589)
$ws.Range("H62").Value = 4575
$ws.Range("I62").Value = 4575
$ws.Range("K62").Value = 4575
$ws.Range("M62").Value = -3951
# Row 65 (anchor G=12589)
$ws.Range("H65").Value = 4575
$ws.Range("I65").Value = 4575
$ws.Range("K65").Value = 22875
$ws.Range("M65").Value = -19755
# Row 136 (anchor G=44031)
$ws.Range("H136").Value = 2858.7144
$ws.Range("I136").Value = 2501.8333
$ws.Range("K136").Value = 7505.499899999999
$ws.Range("M136").Value = -4955.499899999999
